$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data set (player, position, team) in the target row order.
$data = @(
    @("Ja Morant",        "PG",          "Memphis Grizzlies"),
    @("Josh Giddey",      "PG,SG,SF",    "Chicago Bulls"),
    @("Scottie Barnes",   "PG,SG,SF,PF", "Toronto Raptors"),
    @("Mikal Bridges",    "SG,SF,PF",    "New York Knicks"),
    @("Nikola Vucevic",   "PF,C",        "Chicago Bulls"),
    @("KJ Simpson",       "PG",          "Charlotte Hornets"),
    @("P.J. Washington",  "SF,PF",       "Dallas Mavericks"),
    @("Brook Lopez",      "C",           "Milwaukee Bucks"),
    @("Evan Mobley",      "PF,C",        "Cleveland Cavaliers"),
    @("Shaedon Sharpe",   "SG,SF",       "Portland Trail Blazers"),
    @("Luka Doncic",      "PG,SG",       "Los Angeles Lakers"),
    @("De'Aaron Fox",     "PG",          "San Antonio Spurs"),
    @("Isaiah Collier",   "PG,SG",       "Utah Jazz"),
    @("Miles Bridges",    "SF,PF",       "Charlotte Hornets"),
    @("DeMar DeRozan",    "SF,PF",       "Sacramento Kings"),
    @("Tyler Herro",      "PG,SG",       "Miami Heat")
)

# Clear out the previous data range (including the row that will no longer
# be used, since the new table has one fewer row than the old one).
$ws.Range("A2:C18").Clear()

$rowIndex = 2
foreach ($rec in $data) {
    $ws.Cells.Item($rowIndex, 1).Value = $rec[0]
    $ws.Cells.Item($rowIndex, 2).Value = $rec[1]
    $ws.Cells.Item($rowIndex, 3).Value = $rec[2]
    $rowIndex++
}
